$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the editable "N AUD" supplier-number strings in column D with
# plain numeric values (160..320 step 10), removing the old shared strings.
$values = @(160,170,180,190,200,210,220,230,240,250,260,270,280,290,300,310,320)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $values[$i]
}

# New edit button / selection location
$ws.Range("I8").Select()
